$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 8-15 and add new rows 16-17 (columns A:E)
# Row layout: A=id, B=name (string), C=from_bus, D=to_bus, E=in_service (bool)

$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# New rows 16-17 (column A) reuse the same cell formatting as the existing
# "id" cells in column A (bold, centered, thin border) - clone via copy/paste formats
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null
